$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0200434094887143
$ws.Range("D2").Value = 0.01573050031654333
$ws.Range("E2").Value = 2.606204256907688
$ws.Range("F2").Value = 0.3008746626169483
$ws.Range("G2").Value = 0.002347413942568164
$ws.Range("I2").Value = 0.2151339154038396
$ws.Range("M2").Value = 10.73748604104105
$ws.Range("O2").Value = 0.8670485619508099
$ws.Range("C3").Value = 0.01752367395426546
$ws.Range("D3").Value = 0.01376358889438478
$ws.Range("E3").Value = 2.269143849709252
$ws.Range("F3").Value = 0.3044591907494265
$ws.Range("G3").Value = 0.002352826746773228
$ws.Range("I3").Value = 0.2187144283919551
$ws.Range("M3").Value = 9.381662007845648
$ws.Range("O3").Value = 0.8954597981311139
$ws.Range("C4").Value = 0.01596922599944151
$ws.Range("D4").Value = 0.01256053947008695
$ws.Range("E4").Value = 2.062605310405786
$ws.Range("F4").Value = 0.3076306343039761
$ws.Range("G4").Value = 0.002356282461737881
$ws.Range("I4").Value = 0.2216539881105568
$ws.Range("M4").Value = 8.547829241095599
$ws.Range("O4").Value = 0.9163309354213851
$ws.Range("C5").Value = 0.01533396446430402
$ws.Range("D5").Value = 0.01207129960231867
$ws.Range("E5").Value = 1.978516891768066
$ws.Range("F5").Value = 0.3091601265523778
$ws.Range("G5").Value = 0.002357724206886076
$ws.Range("I5").Value = 0.2230328270713287
$ws.Range("M5").Value = 8.207613826726913
$ws.Range("O5").Value = 0.9256720703308332
$ws.Range("C6").Value = 0.0152283713706538
$ws.Range("D6").Value = 0.01199011862145483
$ws.Range("E6").Value = 1.964558017354449
$ws.Range("F6").Value = 0.3094282171745562
$ws.Range("G6").Value = 0.002357965639079443
$ws.Range("I6").Value = 0.2232725542820049
$ws.Range("M6").Value = 8.151093605517929
$ws.Range("O6").Value = 0.9272729082315436
$ws.Range("C7").Value = 0.01596066592879453
$ws.Range("D7").Value = 0.01255393750403755
$ws.Range("E7").Value = 2.061470984222296
$ws.Range("F7").Value = 0.3076503106363973
$ws.Range("G7").Value = 0.0023563017695084
$ws.Range("I7").Value = 0.2216718580606916
$ws.Range("M7").Value = 8.54324278718633
$ws.Range("O7").Value = 0.916453561940699
$ws.Range("C8").Value = 0.0191761385571283
$ws.Range("D8").Value = 0.01505124723004059
$ws.Range("E8").Value = 2.48988170344424
$ws.Range("F8").Value = 0.3019049589907468
$ws.Range("G8").Value = 0.002349252995219228
$ws.Range("I8").Value = 0.2162113719485177
$ws.Range("M8").Value = 10.27022604064683
$ws.Range("O8").Value = 0.8761180054980571
$ws.Range("C9").Value = 0.02542299128452896
$ws.Range("D9").Value = 0.01999332617117489
$ws.Range("E9").Value = 3.334720529646233
$ws.Range("F9").Value = 0.2986571896257857
$ws.Range("G9").Value = 0.002336466885483441
$ws.Range("I9").Value = 0.2116317524458466
$ws.Range("M9").Value = 13.65040484974105
$ws.Range("O9").Value = 0.8253825431057749
$ws.Range("C10").Value = 0.029977023297306
$ws.Range("D10").Value = 0.02366434996396549
$ws.Range("E10").Value = 3.960529526336586
$ws.Range("F10").Value = 0.3016090548561863
$ws.Range("G10").Value = 0.002327685990018582
$ws.Range("I10").Value = 0.2123535080188645
$ws.Range("M10").Value = 16.13675056858943
$ws.Range("O10").Value = 0.8070565006937898
$ws.Range("C11").Value = 0.03204131494655371
$ws.Range("D11").Value = 0.02534628358819901
$ws.Range("E11").Value = 4.246889380416349
$ws.Range("F11").Value = 0.3042120627471618
$ws.Range("G11").Value = 0.00232382006649252
$ws.Range("I11").Value = 0.2136479731319767
$ws.Range("M11").Value = 17.27020394811461
$ws.Range("O11").Value = 0.8032076406487079
$ws.Range("C12").Value = 0.03282197375530416
$ws.Range("D12").Value = 0.02598519517270859
$ws.Range("E12").Value = 4.355617160271379
$ws.Range("F12").Value = 0.3053876228739512
$ws.Range("G12").Value = 0.002322374254861342
$ws.Range("I12").Value = 0.2142838409437644
$ws.Range("M12").Value = 17.69991138302032
$ws.Range("O12").Value = 0.8024276387879468
$ws.Range("C13").Value = 0.03265389099135518
$ws.Range("D13").Value = 0.02584750047082451
$ws.Range("E13").Value = 4.332186992515119
$ws.Range("F13").Value = 0.3051258600359148
$ws.Range("G13").Value = 0.002322684835917062
$ws.Range("I13").Value = 0.214140306030373
$ws.Range("M13").Value = 17.60734176386507
$ws.Range("O13").Value = 0.8025649704699731
$ws.Range("C14").Value = 0.03210556101329587
$ws.Range("D14").Value = 0.02539880545993611
$ws.Range("E14").Value = 4.255828384724907
$ws.Range("F14").Value = 0.3043049211401794
$ws.Range("G14").Value = 0.002323700757504195
$ws.Range("I14").Value = 0.213697327120876
$ws.Range("M14").Value = 17.30554545757985
$ws.Range("O14").Value = 0.8031297033745091
$ws.Range("C15").Value = 0.03176955776010004
$ws.Range("D15").Value = 0.02512423596441948
$ws.Range("E15").Value = 4.209095749459038
$ws.Range("F15").Value = 0.3038270592426571
$ws.Range("G15").Value = 0.002324325388950831
$ws.Range("I15").Value = 0.2134451664364008
$ws.Range("M15").Value = 17.12075535610342
$ws.Range("O15").Value = 0.8035648205804478
$ws.Range("C16").Value = 0.02984196991724275
$ws.Range("D16").Value = 0.02355469571928381
$ws.Range("E16").Value = 3.941853012191928
$ws.Range("F16").Value = 0.301465032347167
$ws.Range("G16").Value = 0.002327941194095173
$ws.Range("I16").Value = 0.2122889245505917
$ws.Range("M16").Value = 16.06273739430964
$ws.Range("O16").Value = 0.8074010588157421
$ws.Range("C17").Value = 0.02865757953205161
$ws.Range("D17").Value = 0.02259508912321451
$ws.Range("E17").Value = 3.778370709816386
$ws.Range("F17").Value = 0.3003448249763352
$ws.Range("G17").Value = 0.002330192051247767
$ws.Range("I17").Value = 0.2118317874197118
$ws.Range("M17").Value = 15.41439389109428
$ws.Range("O17").Value = 0.8109263914984126
$ws.Range("C18").Value = 0.02797565524419099
$ws.Range("D18").Value = 0.02204425692341516
$ws.Range("E18").Value = 3.68449462802198
$ws.Range("F18").Value = 0.2998185916693714
$ws.Range("G18").Value = 0.002331498806981558
$ws.Range("I18").Value = 0.2116593735551078
$ws.Range("M18").Value = 15.04170032476338
$ws.Range("O18").Value = 0.8133744870200985
$ws.Range("C19").Value = 0.02774464786320152
$ws.Range("D19").Value = 0.02185793683683102
$ws.Range("E19").Value = 3.652734947628005
$ws.Range("F19").Value = 0.2996604400906477
$ws.Range("G19").Value = 0.002331943344709952
$ws.Range("I19").Value = 0.2116163421046053
$ws.Range("M19").Value = 14.91554557894841
$ws.Range("O19").Value = 0.8142746167487189
$ws.Range("C20").Value = 0.0287837316864028
$ws.Range("D20").Value = 0.0226971242551528
$ws.Range("E20").Value = 3.795757340467731
$ws.Range("F20").Value = 0.3004517904583892
$ws.Range("G20").Value = 0.002329951191574958
$ws.Range("I20").Value = 0.211871033895612
$ws.Range("M20").Value = 15.48338773190704
$ws.Range("O20").Value = 0.810507373114433
$ws.Range("C21").Value = 0.03226664697666592
$ws.Range("D21").Value = 0.02553054143500333
$ws.Range("E21").Value = 4.278248493353431
$ws.Range("F21").Value = 0.3045408276842849
$ws.Range("G21").Value = 0.002323401867586027
$ws.Range("I21").Value = 0.2138234319608472
$ws.Range("M21").Value = 17.39417571718468
$ws.Range("O21").Value = 0.8029451789552695
$ws.Range("C22").Value = 0.03453686416314383
$ws.Range("D22").Value = 0.02739414795789941
$ws.Range("E22").Value = 4.595297545869698
$ws.Range("F22").Value = 0.3083241241247805
$ws.Range("G22").Value = 0.002319227001168133
$ws.Range("I22").Value = 0.2159518967920704
$ws.Range("M22").Value = 18.64593945661244
$ws.Range("O22").Value = 0.8019675018208261
$ws.Range("C23").Value = 0.03332575445554653
$ws.Range("D23").Value = 0.02639832701758138
$ws.Range("E23").Value = 4.425908578026963
$ws.Range("F23").Value = 0.3062003242695823
$ws.Range("G23").Value = 0.002321445675870604
$ws.Range("I23").Value = 0.2147355986990078
$ws.Range("M23").Value = 17.97752698452081
$ws.Range("O23").Value = 0.8021156028453902
$ws.Range("C24").Value = 0.02872670142153311
$ws.Range("D24").Value = 0.02265099151903627
$ws.Range("E24").Value = 3.787896499285637
$ws.Range("F24").Value = 0.3004030651283927
$ws.Range("G24").Value = 0.002330060044678853
$ws.Range("I24").Value = 0.2118530094950302
$ws.Range("M24").Value = 15.45219547875649
$ws.Range("O24").Value = 0.8106955013803088
$ws.Range("C25").Value = 0.02373934098318387
$ws.Range("D25").Value = 0.01865032090977792
$ws.Range("E25").Value = 3.10544626894719
$ws.Range("F25").Value = 0.2986316687251858
$ws.Range("G25").Value = 0.00233981678878147
$ws.Range("I25").Value = 0.2121789972218977
$ws.Range("M25").Value = 12.73599834161814
$ws.Range("O25").Value = 0.8359030933738438

Write-Output "Updated 192 cells"